$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New diary entry for 2023-08-14 (serial 45152):
# "- Skripteinbindung für automatische Erstellung benötigter DB-Strukturen
#  - Erstellung "Setup-Bereich" für Beladen mit Demodaten per Klick
#  - Prototyping FormBuilder-Klasse für effizienteres Erstellen von Standard-Formularen"
$newText = "- Skripteinbindung für automatische Erstellung benötigter DB-Strukturen`n- Erstellung ""Setup-Bereich"" für Beladen mit Demodaten per Klick`n- Prototyping FormBuilder-Klasse für effizienteres Erstellen von Standard-Formularen"

# A23: date (matches formatting of the other date cells in column A)
$ws.Range("A23").Value = 45152
$ws.Range("A23").NumberFormat = "d-mmm"
$ws.Range("A23").VerticalAlignment = -4160

# B23: activity text, wrapped + quote-prefixed like the other "- ..." bullet entries
$ws.Range("B23").WrapText = $true
$ws.Range("B23").Value = "'" + $newText

# E23: "Umsetzung" hours
$ws.Range("E23").Value = 2.5

# Row height to fit the wrapped text (matches row 3's 3-line height)
$ws.Rows.Item(23).RowHeight = 60

# Leave the same cell selected as in the authored workbook
$ws.Range("F22").Select()
